$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("KPVA Elite Team", $true, $false, $false, $false, $false, $true, 1, $false, "Team Keba", 2)
